$wb = $excel.ActiveWorkbook

# --- Transformers sheet: update row 2 (sub-efficiency tweak + new output product) ---
$wsTrans = $wb.Worksheets.Item("Transformers")
$wsTrans.Range("N2").Value = 0.999
$wsTrans.Range("O2").Value = "km"
$wsTrans.Range("P2").Value = 0.001

# Row 4 and row 5 swap places (product rows reordered)
$wsTrans.Range("A4").Value = "B2gas"
$wsTrans.Range("G4").Value = 0.5
$wsTrans.Range("H4").Value = 0.6
$wsTrans.Range("I4").Value = "biomass"
$wsTrans.Range("M4").Value = "gasoline"

$wsTrans.Range("A5").Value = "Gtkm"
$wsTrans.Range("G5").Value = 0.4
$wsTrans.Range("H5").Value = 0.5
$wsTrans.Range("I5").Value = "gasoline"
$wsTrans.Range("M5").Value = "km"

# --- Connectors sheet: add new connector row (Refinery -> KmHub, in km) ---
$wsConn = $wb.Worksheets.Item("Connectors")
$wsConn.Range("A11").Value = "ref2km"
$wsConn.Range("B11").Value = "Refinery"
$wsConn.Range("C11").Value = "KmHub"
$wsConn.Range("D11").Value = "km"
$wsConn.Range("D11").Select()

# --- Make Transformers the active/selected sheet, with P2 selected ---
$wsTrans.Activate()
$wsTrans.Range("P2").Select()
